$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Day 43428 (2018-11-24), Start 13:00, End 19:39
$ws.Range("A7").Value = 43428
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B7").Value = 0.54166666666666663
$ws.Range("C7").Value = 0.81874999999999998

# Total Hours now exceeds 24h, so switch number format to elapsed-time
$ws.Range("E3").NumberFormat = "[h]:mm:ss"

$wb.Save()
